# Edit script reproducing:
#  1. The table on slide 16 switches its table style from the custom
#     "Table_0" style ({8E8C4BEC-4F72-498A-B004-A550455A2EB6}) to the
#     built-in "No Style, Table Grid" style
#     ({92189ECD-DDDD-467C-8652-D1A23857EF72}).
#  2. The deck's theme (ppt/theme/theme1.xml, used by the slide master)
#     switches its 12 theme colours from the "Integral" palette to the
#     "Office" palette (the rest of the theme - fonts/format scheme -
#     is already identical between the two themes).

$p = $ppt.ActivePresentation

# -- 1. Re-style the table on slide 16 -------------------------------------
$tableShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tableShape = $shape
        }
    }
}

if ($tableShape -ne $null) {
    $tableShape.Table.ApplyStyle("{92189ECD-DDDD-467C-8652-D1A23857EF72}")
}

# -- 2. Re-colour the presentation theme (Integral -> Office) --------------
function Set-ThemeColor($scheme, $index, $rrggbb) {
    $r = [Convert]::ToInt32($rrggbb.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($rrggbb.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($rrggbb.Substring(4, 2), 16)
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$slideOne = $p.Slides.Item(1)
$themeColors = $slideOne.ThemeColorScheme

Set-ThemeColor $themeColors 1  "000000"
Set-ThemeColor $themeColors 2  "FFFFFF"
Set-ThemeColor $themeColors 3  "44546A"
Set-ThemeColor $themeColors 4  "E7E6E6"
Set-ThemeColor $themeColors 5  "5B9BD5"
Set-ThemeColor $themeColors 6  "ED7D31"
Set-ThemeColor $themeColors 7  "A5A5A5"
Set-ThemeColor $themeColors 8  "FFC000"
Set-ThemeColor $themeColors 9  "4472C4"
Set-ThemeColor $themeColors 10 "70AD47"
Set-ThemeColor $themeColors 11 "0563C1"
Set-ThemeColor $themeColors 12 "954F72"
